$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 29
$ws.Range("H29").Value = 7498.2856
$ws.Range("I29").Value = 2251
$ws.Range("J29").Value = 9597.200000000001
$ws.Range("K29").Value = 6753
$ws.Range("L29").Value = 28791.6
$ws.Range("M29").Value = -6472
$ws.Range("N29").Value = -29353.6
# Row 80
$ws.Range("H80").Value = 400.4
$ws.Range("I80").Value = 348.2
$ws.Range("J80").Value = 504.8
$ws.Range("K80").Value = 1044.6
$ws.Range("L80").Value = 1514.4
$ws.Range("M80").Value = -46.59999999999991
$ws.Range("N80").Value = -3510.4
# Row 83
$ws.Range("H83").Value = 400.4
$ws.Range("I83").Value = 348.2
$ws.Range("J83").Value = 504.8
$ws.Range("K83").Value = 3133.8
$ws.Range("L83").Value = 4543.2
$ws.Range("M83").Value = 1858.2
$ws.Range("N83").Value = -14527.2
# Row 88
$ws.Range("H88").Value = 1940.2
$ws.Range("J88").Value = 2627.2856
$ws.Range("L88").Value = 2627.2856
$ws.Range("N88").Value = -3439.2856
# Row 91
$ws.Range("H91").Value = 1940.2
$ws.Range("J91").Value = 2627.2856
$ws.Range("L91").Value = 2627.2856
$ws.Range("N91").Value = -5435.2856
# Row 113
$ws.Range("H113").Value = 3196.625
$ws.Range("I113").Value = 2596.5
$ws.Range("K113").Value = 2596.5
$ws.Range("M113").Value = 657.5
# Row 137
$ws.Range("H137").Value = 2062.3333
$ws.Range("J137").Value = 2430.524
$ws.Range("L137").Value = 7291.572
$ws.Range("N137").Value = -12391.572

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 2726.8333
$ws.Range("I107").Value = 2763.7693
$ws.Range("J107").Value = 2486.75
$ws.Range("K107").Value = 2763.7693
$ws.Range("L107").Value = 2486.75
$ws.Range("M107").Value = -843.7692999999999
$ws.Range("N107").Value = -6326.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 326.94446
$ws.Range("I7").Value = 482.27274
$ws.Range("K7").Value = 482.27274
$ws.Range("M7").Value = -369.27274
# Row 58
$ws.Range("H58").Value = 3066.2856
$ws.Range("I58").Value = 3048.8
$ws.Range("K58").Value = 3048.8
$ws.Range("M58").Value = -2845.8
# Row 99
$ws.Range("H99").Value = 7686.8276
$ws.Range("J99").Value = 8725.727999999999
$ws.Range("L99").Value = 8725.727999999999
$ws.Range("N99").Value = -11721.728
# Row 107
$ws.Range("H107").Value = 15626026
$ws.Range("I107").Value = 25000630
$ws.Range("J107").Value = 1683.6666
$ws.Range("K107").Value = 25000630
$ws.Range("L107").Value = 1683.6666
$ws.Range("M107").Value = -24998710
$ws.Range("N107").Value = -5523.6666
# Row 122
$ws.Range("H122").Value = 2824.5
$ws.Range("I122").Value = 851
$ws.Range("J122").Value = 8745
$ws.Range("K122").Value = 2553
$ws.Range("L122").Value = 26235
$ws.Range("M122").Value = -103
$ws.Range("N122").Value = -31135
# Row 126
$ws.Range("H126").Value = 7686.8276
$ws.Range("J126").Value = 8725.727999999999
$ws.Range("L126").Value = 26177.184
$ws.Range("N126").Value = -31117.184
# Row 136
$ws.Range("H136").Value = 3066.2856
$ws.Range("I136").Value = 3048.8
$ws.Range("K136").Value = 9146.400000000001
$ws.Range("M136").Value = -6596.400000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 21
$ws.Range("H21").Value = 1247.5
$ws.Range("I21").Value = 1247.5
$ws.Range("K21").Value = 3742.5
$ws.Range("M21").Value = -3569.5
# Row 63
$ws.Range("H63").Value = 5433.3335
$ws.Range("I63").Value = 5433.3335
$ws.Range("K63").Value = 16300.0005
$ws.Range("M63").Value = -15551.0005
# Row 66
$ws.Range("H66").Value = 5433.3335
$ws.Range("I66").Value = 5433.3335
$ws.Range("K66").Value = 48900.0015
$ws.Range("M66").Value = -45156.0015
# Row 68
$ws.Range("H68").Value = 1367.25
$ws.Range("J68").Value = 1367.25
$ws.Range("L68").Value = 4101.75
$ws.Range("N68").Value = -5723.75
# Row 69
$ws.Range("H69").Value = 897.8333
$ws.Range("I69").Value = 897.8333
$ws.Range("K69").Value = 2693.4999
$ws.Range("M69").Value = -1882.4999
# Row 70 (N70 is removed entirely; M70 becomes the new last value)
$ws.Range("H70").Value = 5972.2856
$ws.Range("I70").Value = 5972.2856
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 17916.8568
$ws.Range("L70").Value = 0
$ws.Range("M70").Value = -17601.8568
$ws.Range("N70").ClearContents()
# Row 71
$ws.Range("H71").Value = 1367.25
$ws.Range("J71").Value = 1367.25
$ws.Range("L71").Value = 12305.25
$ws.Range("N71").Value = -20417.25
# Row 72
$ws.Range("H72").Value = 897.8333
$ws.Range("I72").Value = 897.8333
$ws.Range("K72").Value = 8080.4997
$ws.Range("M72").Value = -4024.4997
# Row 73 (N73 is removed entirely; M73 becomes the new last value)
$ws.Range("H73").Value = 5972.2856
$ws.Range("I73").Value = 5972.2856
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 17916.8568
$ws.Range("L73").Value = 0
$ws.Range("M73").Value = -16824.8568
$ws.Range("N73").ClearContents()
# Row 140
$ws.Range("H140").Value = 2675.1853
$ws.Range("I140").Value = 1405.7826
$ws.Range("K140").Value = 4217.3478
$ws.Range("M140").Value = 962.6522000000004

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3317.4375
$ws.Range("I132").Value = 3215.3057
$ws.Range("K132").Value = 9645.917099999999
$ws.Range("M132").Value = -7115.917099999999

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2272.7778
$ws.Range("J22").Value = 2874.75
$ws.Range("L22").Value = 2874.75
$ws.Range("N22").Value = -3464.75
# Row 27
$ws.Range("H27").Value = 2272.7778
$ws.Range("J27").Value = 2874.75
$ws.Range("L27").Value = 2874.75
$ws.Range("N27").Value = -3088.75
# Row 61
$ws.Range("H61").Value = 2351.5386
$ws.Range("I61").Value = 2157.3
$ws.Range("K61").Value = 2157.3
$ws.Range("M61").Value = -1955.3
# Row 82
$ws.Range("H82").Value = 2838.7058
$ws.Range("I82").Value = 2957.5833
$ws.Range("K82").Value = 2957.5833
$ws.Range("M82").Value = -2596.5833
# Row 85
$ws.Range("H85").Value = 2838.7058
$ws.Range("I85").Value = 2957.5833
$ws.Range("K85").Value = 2957.5833
$ws.Range("M85").Value = -1709.5833
# Row 113
$ws.Range("H113").Value = 2351.5386
$ws.Range("I113").Value = 2157.3
$ws.Range("K113").Value = 2157.3
$ws.Range("M113").Value = 12.69999999999982
# Row 122
$ws.Range("H122").Value = 3663.1428
$ws.Range("I122").Value = 2935
$ws.Range("J122").Value = 6333
$ws.Range("K122").Value = 8805
$ws.Range("L122").Value = 18999
$ws.Range("M122").Value = -6355
$ws.Range("N122").Value = -23899
# Row 136
$ws.Range("H136").Value = 5091.522
$ws.Range("I136").Value = 4142.7646
$ws.Range("J136").Value = 7779.6665
$ws.Range("K136").Value = 12428.2938
$ws.Range("L136").Value = 23338.9995
$ws.Range("M136").Value = -9878.293800000001
$ws.Range("N136").Value = -28438.9995

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 431.06668
$ws.Range("I107").Value = 303
$ws.Range("K107").Value = 909
$ws.Range("M107").Value = 1011
